$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New song rows appended below the existing data (row 23 and row 34 are
# intentionally left blank, matching the gaps already present in the sheet).
$rows = @{
    24 = @("A2", "David Bowie", "Space Oddity", "/music/David Bowie/Bowie_ The Singles 1969-1993 (Disc 1)/01 Space Oddity.wav")
    25 = @("B2", "David Bowie", "Life On Mars", "/music/David Bowie/Bowie_ The Singles 1969-1993 (Disc 1)/04 Life on Mars_.wav")
    26 = @("C2", "David Bowie", "Starman", "/music/David Bowie/Bowie_ The Singles 1969-1993 (Disc 1)/06 Starman.wav")
    27 = @("D2", "David Bowie", "Fame", "/music/David Bowie/Bowie_ The Singles 1969-1993 (Disc 1)/15 Fame.wav")
    28 = @("E2", "David Bowie", "Under Pressure", "/music/David Bowie/Bowie_ The Singles 1969-1993 (Disc 2)/08 Under Pressure [with Queen].wav ")
    29 = @("F2", "David Bowie", "Dancing in the Street", "/music/David Bowie/Bowie_ The Singles 1969-1993 (Disc 2)/15 Dancing in the Street [with Mick Jagger].wav")
    30 = @("G2", "Dire Straits", "Money For Nothing", "/music/Dire Straits/Brothers In Arms/02 Money for Nothing.wav")
    31 = @("H2", "Dire Straits", "On Every Street", "/music/Dire Straits/On Every Street/02 On Every Street.wav")
    32 = @("J2", "Duran Duran", "A View To A Kill", "/music/Duran Duran/Greatest/03 A View To A Kill.wav  ")
    33 = @("K2", "Duran Duran", "Hungry Like the Wolf", "/music/Duran Duran/Greatest/07 Hungry Like The Wolf.wav")
    35 = @("L2", "Elton John", "I'm Still Standing", "/music/Elton John/Greatest Hits 1976-1986/01 I'm Still Standing.wav ")
    36 = @("M2", "Elton John", "Don't Go Breaking My Heart", "/music/Elton John/Greatest Hits 1976-1986/06 Don't Go Breaking My Heart.wav ")
}

foreach ($r in ($rows.Keys | Sort-Object)) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
}

# Widen column D to fit the new, longer filenames.
$ws.Columns.Item(4).ColumnWidth = 109.42

# Scroll/selection state, matching where editing left off.
$ws.Activate() | Out-Null
$ws.Range("D33").Select() | Out-Null
